$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H137").Value = 19608966
$ws.Range("I137").Value = 24391026
$ws.Range("K137").Value = 73173078
$ws.Range("M137").Value = -73170528
$ws.Range("H138").Value = 3743606.5
$ws.Range("I138").Value = 1325216
$ws.Range("J138").Value = 4697338
$ws.Range("K138").Value = 3975648
$ws.Range("L138").Value = 14092014
$ws.Range("M138").Value = -3970508
$ws.Range("N138").Value = -14102294

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H22").Value = 573.2
$ws.Range("I22").Value = 573.2
$ws.Range("K22").Value = 573.2
$ws.Range("M22").Value = -274.2
$ws.Range("H61").Value = 2874.7407
$ws.Range("I61").Value = 1748
$ws.Range("J61").Value = 5128.222
$ws.Range("K61").Value = 1748
$ws.Range("L61").Value = 5128.222
$ws.Range("M61").Value = -1536
$ws.Range("N61").Value = -5552.222
$ws.Range("H74").Value = 5161.8438
$ws.Range("I74").Value = 1294.2106
$ws.Range("K74").Value = 1294.2106
$ws.Range("M74").Value = -420.2106000000001
$ws.Range("H77").Value = 5161.8438
$ws.Range("I77").Value = 1294.2106
$ws.Range("K77").Value = 6471.053000000001
$ws.Range("M77").Value = -2103.053000000001
$ws.Range("H132").Value = 3999.4167
$ws.Range("I132").Value = 3832.4
$ws.Range("J132").Value = 4834.5
$ws.Range("K132").Value = 11497.2
$ws.Range("L132").Value = 14503.5
$ws.Range("M132").Value = -8967.200000000001
$ws.Range("N132").Value = -19563.5
$ws.Range("H133").Value = 33614.145
$ws.Range("J133").Value = 33614.145
$ws.Range("L133").Value = 33614.145
$ws.Range("N133").Value = -38674.145
$ws.Range("H136").Value = 2874.7407
$ws.Range("I136").Value = 1748
$ws.Range("J136").Value = 5128.222
$ws.Range("K136").Value = 5244
$ws.Range("L136").Value = 15384.666
$ws.Range("M136").Value = -2694
$ws.Range("N136").Value = -20484.666

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H134").Value = 2848.6562
$ws.Range("I134").Value = 2189
$ws.Range("K134").Value = 6567
$ws.Range("M134").Value = -4032
$ws.Range("H135").Value = 43870
$ws.Range("J135").Value = 43870
$ws.Range("L135").Value = 43870
$ws.Range("N135").Value = -54010

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 956.1053000000001
$ws.Range("I31").Value = 903.3714
$ws.Range("J31").Value = 1571.3334
$ws.Range("K31").Value = 903.3714
$ws.Range("L31").Value = 1571.3334
$ws.Range("M31").Value = -608.3714
$ws.Range("N31").Value = -2161.3334
$ws.Range("H34").Value = 956.1053000000001
$ws.Range("I34").Value = 903.3714
$ws.Range("J34").Value = 1571.3334
$ws.Range("K34").Value = 903.3714
$ws.Range("L34").Value = 1571.3334
$ws.Range("M34").Value = -701.3714
$ws.Range("N34").Value = -1975.3334
$ws.Range("H58").Value = 1899.3846
$ws.Range("J58").Value = 4173.2
$ws.Range("L58").Value = 4173.2
$ws.Range("N58").Value = -4579.2
$ws.Range("H132").Value = 2453.2083
$ws.Range("I132").Value = 2139.9092
$ws.Range("J132").Value = 5899.5
$ws.Range("K132").Value = 6419.7276
$ws.Range("L132").Value = 17698.5
$ws.Range("M132").Value = -3889.7276
$ws.Range("N132").Value = -22758.5
$ws.Range("H134").Value = 2072.4443
$ws.Range("I134").Value = 1095.138
$ws.Range("J134").Value = 6121.2856
$ws.Range("K134").Value = 3285.414
$ws.Range("L134").Value = 18363.8568
$ws.Range("M134").Value = -750.4139999999998
$ws.Range("N134").Value = -23433.8568
$ws.Range("H136").Value = 1899.3846
$ws.Range("J136").Value = 4173.2
$ws.Range("L136").Value = 12519.6
$ws.Range("N136").Value = -17619.6
$ws.Range("H140").Value = 33853.332
$ws.Range("J140").Value = 45780
$ws.Range("L140").Value = 45780
$ws.Range("N140").Value = -56140

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H64").Value = 2983.0908
$ws.Range("J64").Value = 3181.4
$ws.Range("L64").Value = 9544.200000000001
$ws.Range("N64").Value = -10084.2
$ws.Range("H67").Value = 2983.0908
$ws.Range("J67").Value = 3181.4
$ws.Range("L67").Value = 9544.200000000001
$ws.Range("N67").Value = -11416.2

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H80").Value = 2319.2222
$ws.Range("I80").Value = 2244
$ws.Range("K80").Value = 2244
$ws.Range("M80").Value = -1246
$ws.Range("H83").Value = 2319.2222
$ws.Range("I83").Value = 2244
$ws.Range("K83").Value = 11220
$ws.Range("M83").Value = -6228
$ws.Range("H126").Value = 2196.6572
$ws.Range("I126").Value = 1681.8125
$ws.Range("K126").Value = 5045.4375
$ws.Range("M126").Value = -2575.4375
$ws.Range("H132").Value = 3262.2
$ws.Range("I132").Value = 3040.5625
$ws.Range("K132").Value = 9121.6875
$ws.Range("M132").Value = -6591.6875

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H82").Value = 1294.0625
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1294.0625
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 1294.0625
$ws.Range("M82").Value = $null
$ws.Range("N82").Value = -2016.0625
$ws.Range("H85").Value = 1294.0625
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1294.0625
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 1294.0625
$ws.Range("M85").Value = $null
$ws.Range("N85").Value = -3790.0625
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").Value = $null
$ws.Range("H132").Value = 4659.2563
$ws.Range("I132").Value = 3839.32
$ws.Range("J132").Value = 6123.4287
$ws.Range("K132").Value = 11517.96
$ws.Range("L132").Value = 18370.2861
$ws.Range("M132").Value = -8987.960000000001
$ws.Range("N132").Value = -23430.2861
$ws.Range("H133").Value = 44931.4
$ws.Range("J133").Value = 44931.4
$ws.Range("L133").Value = 44931.4
$ws.Range("N133").Value = -49991.4
$ws.Range("H136").Value = 3644.3865
$ws.Range("I136").Value = 2131.2812
$ws.Range("J136").Value = 7679.3335
$ws.Range("K136").Value = 6393.8436
$ws.Range("L136").Value = 23038.0005
$ws.Range("M136").Value = -3843.8436
$ws.Range("N136").Value = -28138.0005

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H32").Value = 6842
$ws.Range("I32").Value = 4263
$ws.Range("J32").Value = 12000
$ws.Range("K32").Value = 4263
$ws.Range("L32").Value = 12000
$ws.Range("M32").Value = -3946
$ws.Range("N32").Value = -12634
$ws.Range("H107").Value = 489
$ws.Range("I107").Value = 528.6667
$ws.Range("J107").Value = 370
$ws.Range("K107").Value = 1586.0001
$ws.Range("L107").Value = 1110
$ws.Range("M107").Value = 333.9999
$ws.Range("N107").Value = -4950
$ws.Range("H122").Value = 69467.2
$ws.Range("I122").Value = 127088.5
$ws.Range("K122").Value = 381265.5
$ws.Range("M122").Value = -378815.5
$ws.Range("H132").Value = 15627797
$ws.Range("I132").Value = 16131878
$ws.Range("J132").Value = 1305
$ws.Range("K132").Value = 48395634
$ws.Range("L132").Value = 3915
$ws.Range("M132").Value = -48393104
$ws.Range("N132").Value = -8975
$ws.Range("H136").Value = 11941948
$ws.Range("I136").Value = 16717702
$ws.Range("J136").Value = 2561.875
$ws.Range("K136").Value = 50153106
$ws.Range("L136").Value = 7685.625
$ws.Range("M136").Value = -50150556
$ws.Range("N136").Value = -12785.625
